# Update odds/values on Sheet1 to match the 2024-10-31 FlashScore refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 2 (AS Roma - Torino) ---
$ws.Range("N2").Value = 10
$ws.Range("Q2").Value = 2.05
$ws.Range("R2").Value = 1.68
$ws.Range("AH2").Value = 10
$ws.Range("AT2").Value = 2.62
$ws.Range("BA2").Value = 101

# --- Row 3 (Como - Lazio) ---
$ws.Range("G3").Value = 3.1
$ws.Range("I3").Value = 2.4
$ws.Range("J3").Value = 3.5
$ws.Range("L3").Value = 3.1
$ws.Range("N3").Value = 9.5
$ws.Range("Q3").Value = 2.05
$ws.Range("R3").Value = 1.85
$ws.Range("W3").Value = 9.5
$ws.Range("Y3").Value = 11
$ws.Range("Z3").Value = 29
$ws.Range("AA3").Value = 23
$ws.Range("AH3").Value = 8.5
$ws.Range("AI3").Value = 12
$ws.Range("AK3").Value = 23
$ws.Range("AR3").Value = 67
$ws.Range("AW3").Value = 4.5
$ws.Range("BC3").Value = 126

# --- Row 7 (Grasshoppers - Lugano) ---
$ws.Range("N7").Value = 13

# --- Row 8 (Servette - Luzern) ---
$ws.Range("G8").Value = 1.76
$ws.Range("I8").Value = 4.33
$ws.Range("J8").Value = 2.37
$ws.Range("K8").Value = 2.37
$ws.Range("L8").Value = 4.33
$ws.Range("AE8").Value = 13
$ws.Range("AJ8").Value = 15
$ws.Range("AM8").Value = 34
$ws.Range("AO8").Value = 9
$ws.Range("BA8").Value = 81
$ws.Range("BB8").Value = 151
